# Update cryptocurrency price/volume figures as refreshed by the
# "Updated cryptos list" GitHub Actions workflow.
#
# Prices/volumes in this sheet are stored as plain text (not numbers),
# so every write below forces the cell to Text format first (preserving
# its original Style afterwards) to stop Excel's automatic "looks like a
# number" conversion from silently turning e.g. "580.46" into a numeric
# cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$updates = @(
    @{ Row = 2;  D = "67.775.91";   E = "  +1.09%  " },
    @{ Row = 3;  D = "3.337.26";    E = "  +1.85%  " },
    @{ Row = 4;  D = $null;         E = "  -0.02%  " },
    @{ Row = 5;  D = "580.46";      E = "  +0.66%  " },
    @{ Row = 6;  D = "176.08";      E = "  +2.50%  " },
    @{ Row = 7;  D = $null;         E = "  -0.01%  " },
    @{ Row = 8;  D = $null;         E = "  +2.04%  " },
    @{ Row = 9;  D = "3.332.38";    E = "  +1.88%  " },
    @{ Row = 10; D = "0.183";       E = "  +6.47%  " },
    @{ Row = 11; D = "0.581";       E = "  +2.04%  " },
    @{ Row = 12; D = "46.92";       E = "  +4.67%  " },
    @{ Row = 13; D = $null;         E = "  +1.25%  " },
    @{ Row = 14; D = "690.15";      E = "  +0.34%  " },
    @{ Row = 15; D = "3.877.68";    E = $null },
    @{ Row = 16; D = $null;         E = "  +2.70%  " },
    @{ Row = 17; D = "67.841.87";   E = "  +0.99%  " },
    @{ Row = 19; D = "3.343.11";    E = "  +2.18%  " },
    @{ Row = 20; D = "17.58";       E = "  +2.49%  " },
    @{ Row = 21; D = "11.06";       E = "  +4.20%  " },
    @{ Row = 22; D = "0.894";       E = "  +1.65%  " },
    @{ Row = 23; D = $null;         E = "  +3.79%  " },
    @{ Row = 24; D = "16.97";       E = "  +0.68%  " },
    @{ Row = 25; D = "101.09";      E = "  +3.17%  " },
    @{ Row = 26; D = "3.91";        E = $null },
    @{ Row = 27; D = $null;         E = "  +2.44%  " },
    @{ Row = 28; D = $null;         E = "  +5.87%  " },
    @{ Row = 29; D = "33.05";       E = "  -0.53%  " },
    @{ Row = 30; D = $null;         E = "  +3.27%  " },
    @{ Row = 31; D = $null;         E = "  +7.30%  " },
    @{ Row = 32; D = "568.86";      E = "  -2.24%  " },
    @{ Row = 33; D = "11.01";       E = "  +2.15%  " },
    @{ Row = 34; D = $null;         E = "  +3.30%  " },
    @{ Row = 35; D = "57.45";       E = "  +3.80%  " },
    @{ Row = 36; D = $null;         E = "  -0.16%  " },
    @{ Row = 37; D = "3.705.07";    E = "  -2.77%  " },
    @{ Row = 38; D = $null;         E = "  +1.87%  " },
    @{ Row = 39; D = "35.31";       E = "  +12.96%  " },
    @{ Row = 40; D = $null;         E = "  +4.97%  " },
    @{ Row = 41; D = "3.18";        E = "  +7.43%  " },
    @{ Row = 42; D = $null;         E = "  +3.05%  " },
    @{ Row = 43; D = "0.0₃0674";    E = "  +3.21%  " },
    @{ Row = 44; D = $null;         E = "  +3.94%  " },
    @{ Row = 45; D = "3.31";        E = "  -2.32%  " },
    @{ Row = 46; D = $null;         E = "  +1.97%  " },
    @{ Row = 47; D = $null;         E = "  +5.75%  " },
    @{ Row = 48; D = $null;         E = "  +1.98%  " },
    @{ Row = 49; D = $null;         E = "  -0.12%  " },
    @{ Row = 50; D = $null;         E = "  -0.39%  " },
    @{ Row = 51; D = "132.08";      E = "  +2.80%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-TextValue $u.Row 4 $u.D
    }
    if ($null -ne $u.E) {
        Set-TextValue $u.Row 5 $u.E
    }
}
